$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.910.74"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "3.507.41"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'607.48"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").Value = "'198.38"
$ws.Range("E6").Value = "  +6.55%  "
$ws.Range("E7").Value = "  +1.50%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -1.96%  "
$ws.Range("D10").Value = "'0.658"
$ws.Range("E10").Value = "  +1.92%  "
$ws.Range("D11").Value = "'54.44"
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("D13").Value = "'9.62"
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("D14").Value = "4.068.87"
$ws.Range("E14").Value = "  -1.25%  "
$ws.Range("D15").Value = "'598.34"
$ws.Range("E15").Value = "  +4.56%  "
$ws.Range("D16").Value = "69.983.85"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").Value = "'12.71"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").Value = "3.504.93"
$ws.Range("E19").Value = "  -1.08%  "
$ws.Range("D21").Value = "'0.996"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("E22").Value = "  +2.13%  "
$ws.Range("D23").Value = "'103.77"
$ws.Range("E23").Value = "  +10.20%  "
$ws.Range("D24").Value = "'4.66"
$ws.Range("E24").Value = "  -1.46%  "
$ws.Range("D25").Value = "'5.08"
$ws.Range("E25").Value = "  +4.17%  "
$ws.Range("D26").Value = "'3.12"
$ws.Range("E26").Value = "  +5.91%  "
$ws.Range("D27").Value = "'11.00"
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("D28").Value = "'9.85"
$ws.Range("E28").Value = "  +4.99%  "
$ws.Range("D29").Value = "'33.76"
$ws.Range("E29").Value = "  +4.62%  "
$ws.Range("D30").Value = "'4.61"
$ws.Range("E30").Value = "  +23.72%  "
$ws.Range("E31").Value = "  +2.81%  "
$ws.Range("D32").Value = "'12.79"
$ws.Range("E32").Value = "  +4.54%  "
$ws.Range("E33").Value = "  +1.24%  "
$ws.Range("D34").Value = "'63.77"
$ws.Range("E34").Value = "  -0.41%  "
$ws.Range("D35").Value = "3.711.55"
$ws.Range("E35").Value = "  +3.57%  "
$ws.Range("D36").Value = "'523.70"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").Value = "0.0₃0806"
$ws.Range("E38").Value = "  +3.10%  "
$ws.Range("D39").Value = "'3.01"
$ws.Range("E39").Value = "  -5.21%  "
$ws.Range("D40").Value = "'0.394"
$ws.Range("E40").Value = "  -2.68%  "
$ws.Range("D41").Value = "'36.90"
$ws.Range("E41").Value = "  -1.41%  "
$ws.Range("D42").Value = "'3.56"
$ws.Range("E42").Value = "  +0.66%  "
$ws.Range("D43").Value = "'0.138"
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("D44").Value = "'0.0460"
$ws.Range("E44").Value = "  +0.73%  "
$ws.Range("D45").Value = "'2.87"
$ws.Range("E45").Value = "  -1.91%  "
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("E47").Value = "  -4.09%  "
$ws.Range("D48").Value = "'8.78"
$ws.Range("E48").Value = "  -4.22%  "
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("D50").Value = "'131.98"
$ws.Range("E50").Value = "  -3.44%  "
$ws.Range("D51").Value = "'0.000241"
$ws.Range("E51").Value = "  -1.51%  "
